$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 116, shifting existing rows 116-120 down to 117-121.
$ws.Rows.Item(116).Insert()

# New row 116 retains most fields from the (old) row 116 / now row 117, but with
# updated Fecha, Volumen, Precio minimo/maximo/promedio, Origen and Precio $/Kg.
$ws.Cells.Item(116, 1).Value = 4
$ws.Cells.Item(116, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(116, 3).Value = "Los Lagos"
$ws.Cells.Item(116, 4).Value = 44747
$ws.Cells.Item(116, 5).Value = 10
$ws.Cells.Item(116, 6).Value = 100112022
$ws.Cells.Item(116, 7).Value = "Arveja Verde"
$ws.Cells.Item(116, 8).Value = "Perfection"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 40
$ws.Cells.Item(116, 11).Value = 44000
$ws.Cells.Item(116, 12).Value = 44000
$ws.Cells.Item(116, 13).Value = 44000
$ws.Cells.Item(116, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(116, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(116, 16).Value = 1760
$ws.Cells.Item(116, 17).Value = 25
$ws.Cells.Item(116, 18).Value = "Hortaliza"
